$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to Text format first,
# otherwise Excel auto-converts them to numbers and mangles trailing zeros
# (e.g. "3.600" -> 3.6, "0.001090" -> 1.09E-3). All these source cells were
# originally stored as inline text, so we preserve that by pre-formatting.
$numericLookingCells = @("D2", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D50")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '264.47'
$ws.Range("D4").Value = '6.287'
$ws.Range("D5").Value = '0.06154'
$ws.Range("D6").Value = '3.599'
$ws.Range("D7").Value = '6.681'
$ws.Range("D8").Value = '1.348'
$ws.Range("D9").Value = '0.8302'
$ws.Range("D10").Value = '0.01354'
$ws.Range("D11").Value = '0.1588'
$ws.Range("D12").Value = '0.08239'
$ws.Range("B15").Value = 'ProBitToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D15").Value = '0.1227'
$ws.Range("E15").Value = '14ProBitTokenPROB'
$ws.Range("B16").Value = 'BitMartToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D16").Value = '0.09256'
$ws.Range("E16").Value = '15BitMartTokenBMX'
$ws.Range("B17").Value = 'MCDex'
$ws.Range("C17").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D17").Value = '3.941'
$ws.Range("E17").Value = '16MCDexMCB'
$ws.Range("B18").Value = 'BitForexToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D18").Value = '0.001715'
$ws.Range("E18").Value = '17BitForexTokenBF'
$ws.Range("B19").Value = 'CoinExToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D19").Value = '0.04879'
$ws.Range("E19").Value = '18CoinExTokenCET'
$ws.Range("B20").Value = 'TigerCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D20").Value = '0.006258'
$ws.Range("E20").Value = '19TigerCashTCH'
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").Value = '0.005274'
$ws.Range("E21").Value = '20HotbitTokenHTB'
$ws.Range("B22").Value = 'BitKan'
$ws.Range("C22").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D22").Value = '0.001090'
$ws.Range("E22").Value = '21BitKanKAN'
$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").Value = '0.0001501'
$ws.Range("E23").Value = '22NitroExNTX'
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").Value = '3.766'
$ws.Range("E24").Value = '23LEOLEO'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").Value = '2.289'
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").Value = '0.3378'
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'
$ws.Range("D27").Value = '0.0002679'
$ws.Range("E27").Value = '26UpBotsUBXT'
$ws.Range("D40").Value = '0.04615'
$ws.Range("D41").Value = '0.006951'
$ws.Range("D42").Value = '0.1136'
$ws.Range("D43").Value = '0.003401'
$ws.Range("D44").Value = '0.01074'
$ws.Range("E44").Value = '43LocalTradersLCTWorstin24h'
$ws.Range("D45").Value = '0.00006156'
$ws.Range("D47").Value = '0.7780'
$ws.Range("D48").Value = '0.1956'
$ws.Range("D50").Value = '0.01240'
